$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Part 1: " branch custom name " -> " branch dev (custom name) "
#          with "dev" and "custom name" colored red (bold retained)
# ------------------------------------------------------------------

$rng = $d.Content
$rng.Find.Execute("branch custom name ", $true, $false, $false, $false, $false, $true, 1, $false, "branch dev (custom name) ", 2)

$p21 = $d.Paragraphs.Item(21).Range
$p21Start = $p21.Start
$p21End = $p21.End

# Color "dev"
$f1 = $d.Range($p21Start, $p21End)
$f1.Find.Execute("dev", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$f1.Font.Color = 255

# Color "custom name"
$f2 = $d.Range($p21Start, $p21End)
$f2.Find.Execute("custom name", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$f2.Font.Color = 255

# ------------------------------------------------------------------
# Part 2: append a new list paragraph describing the checkout steps
# ------------------------------------------------------------------

$p21 = $d.Paragraphs.Item(21).Range
$pEnd = $p21.End
# insertion point right after "for create demo branch" (non-bold run) -
# new text typed here naturally inherits the plain (non-bold) formatting.
$tailPoint = $d.Range($pEnd - 1, $pEnd - 1)
$markStart = $tailPoint.Start
$fullText = "git checkout dev change master branch to dev branch."
$tailPoint.InsertAfter($fullText)

$scopeEnd = $d.Paragraphs.Item(21).Range.End
$cursor = $markStart

$s1 = $d.Range($cursor, $scopeEnd)
$s1.Find.Execute("git checkout dev ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s1.Bold = 1
$cursor = $s1.End

$s2 = $d.Range($cursor, $scopeEnd)
$s2.Find.Execute("master ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s2.Bold = 1
$cursor = $s2.End

$s3 = $d.Range($cursor, $scopeEnd)
$s3.Find.Execute("dev ", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$s3.Bold = 1
$cursor = $s3.End

# split the appended text off into its own (already-correctly-styled) paragraph
$splitPoint = $d.Range($markStart, $markStart)
$splitPoint.InsertParagraphBefore()
